$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sexo" and "edad-grupos-quinquenales-2010" columns are re-classified
# from curated dimensions to curated measures. Update the per-column
# metadata rows (type tag, dim/medida flag, data-type) and drop the old
# mapping-file row that only applied to dimension columns.

$ws.Range("A2").Value = "iaest-measure:edad-grupos-quinquenales-2010"
$ws.Range("F2").Value = "iaest-measure:sexo"

$ws.Range("A3").Value = "medida"
$ws.Range("F3").Value = "medida"

$ws.Range("A4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Row 5 (mapping-*.xlsx references) no longer applies now that both
# columns are measures rather than curated dimensions.
$ws.Range("A5:I5").Delete()
